# AltaClienteControlDual - "update entregable 1 y 2"
# The approval workflow re-ran: the "Fecha" (date/time) stamp on row 2
# advances from the April timestamp to the new June timestamp, and the
# active selection moves from K10 to N7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fecha column (P2) with the new approval timestamp.
$ws.Range("P2").Value = "26 jun. 2023, 18:16:32"

# Restore the other workflow columns to their (unchanged) text values,
# keeping them as text rather than letting automatic type-inference turn
# the numeric-looking ones into numbers.
$ws.Range("C2").Value = "D75356820"
$ws.Range("D2").Value = "'75356819"
$ws.Range("M2").Value = "ASARMIENTOS1"

# Move the active selection to where the user last clicked.
$ws.Range("N7").Select() | Out-Null
